# Daily attendance processing - 2025-11-22 22:46:08
# Swap the order of "dnasr281@gmail.com" and "System" in column G
# (Recorded By / Attendance Taker) wherever the old combined value is found.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Session Analysis Results")

$oldValue = "dnasr281@gmail.com, System"
$newValue = "System, dnasr281@gmail.com"

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count + $usedRange.Row - 1

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)  # Column G is the 7th column
    if ($cell.Value2 -eq $oldValue) {
        $cell.Value = $newValue
    }
}
